# Update "想去人数" (column F) figures on the "展览", "演出" and "全部类型"
# sheets to match the newly generated data output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value2  = 6808
$ws1.Range("F3").Value2  = 89
$ws1.Range("F4").Value2  = 21
$ws1.Range("F6").Value2  = 146
$ws1.Range("F7").Value2  = 6443
$ws1.Range("F8").Value2  = 57
$ws1.Range("F9").Value2  = 197
$ws1.Range("F10").Value2 = 1279
$ws1.Range("F11").Value2 = 11
$ws1.Range("F13").Value2 = 395
$ws1.Range("F15").Value2 = 18
$ws1.Range("F16").Value2 = 374
$ws1.Range("F18").Value2 = 8
$ws1.Range("F19").Value2 = 4791
$ws1.Range("F20").Value2 = 84
$ws1.Range("F21").Value2 = 63
$ws1.Range("F22").Value2 = 257
$ws1.Range("F24").Value2 = 134

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value2 = 44

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value2  = 6808
$ws4.Range("F3").Value2  = 89
$ws4.Range("F4").Value2  = 21
$ws4.Range("F7").Value2  = 6443
$ws4.Range("F8").Value2  = 57
$ws4.Range("F9").Value2  = 197
$ws4.Range("F10").Value2 = 1279
$ws4.Range("F11").Value2 = 11
$ws4.Range("F13").Value2 = 395
$ws4.Range("F16").Value2 = 374
$ws4.Range("F18").Value2 = 8
$ws4.Range("F19").Value2 = 4791
$ws4.Range("F20").Value2 = 44
$ws4.Range("F21").Value2 = 84
$ws4.Range("F22").Value2 = 63
$ws4.Range("F23").Value2 = 257
$ws4.Range("F25").Value2 = 134

$wb.Save()
